$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force Text format, assign, then restore the default "Normal" style so the
# cell keeps no extra formatting (matches original unstyled inline-string cells).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.612"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0883"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0357"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "82.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "114.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.29"
$ws.Range("D50").Style = "Normal"

# Remaining cells are unambiguous text already (percent strings with padding
# spaces, or multi-dot price strings) and need no special handling.
$ws.Range("D2").Value = "42.771.66"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "2.300.73"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "2.651.33"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "2.284.19"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "42.719.55"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("E21").Value = "  +35.75%  "
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("E24").Value = "  -6.75%  "
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  +3.02%  "
$ws.Range("E29").Value = "  +9.68%  "
$ws.Range("E30").Value = "  -3.38%  "
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("E32").Value = "  +5.04%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("E35").Value = "  -1.67%  "
$ws.Range("E36").Value = "  -13.87%  "
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("E39").Value = "  +3.35%  "
$ws.Range("E40").Value = "  -3.04%  "
$ws.Range("E41").Value = "  +4.32%  "
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("E45").Value = "  -9.10%  "
$ws.Range("E46").Value = "  +3.02%  "
$ws.Range("E47").Value = "  +4.97%  "
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("D51").Value = "1.602.71"
$ws.Range("E51").Value = "  +2.81%  "
